$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on Column D cells so numeric-looking strings
# (e.g. "1.00", "322.29") are stored as text, matching the original inlineStr cells,
# then restore the default "Normal" style so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.450.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.485.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.26%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.878.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.480.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.324.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0934"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +5.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.73%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.98%  "
$ws.Range("E31").Value = "  +6.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0780"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.14%  "
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "121.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.959.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.61%  "
